$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values are treated as text so multi-dot
# "thousand-separated" price strings and single-dot values keep their
# exact original textual formatting (e.g. trailing zeros) instead of
# being auto-converted to numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '43.058.41'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '2.301.89'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '302.77'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').Value = '99.22'
$ws.Range('E6').Value = '  +5.58%  '
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +3.11%  '
$ws.Range('D10').Value = '34.40'
$ws.Range('D11').Value = '0.0798'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').Value = '49.29'
$ws.Range('E12').Value = '  +3.86%  '
$ws.Range('D13').Value = '0.117'
$ws.Range('E13').Value = '  +4.40%  '
$ws.Range('E14').Value = '  +17.43%  '
$ws.Range('E15').Value = '  +1.96%  '
$ws.Range('D16').Value = '2.662.16'
$ws.Range('E16').Value = '  +1.97%  '
$ws.Range('D17').Value = '2.281.63'
$ws.Range('E17').Value = '  +0.95%  '
$ws.Range('E18').Value = '  +4.44%  '
$ws.Range('D19').Value = '42.965.28'
$ws.Range('E19').Value = '  +2.13%  '
$ws.Range('D20').Value = '12.35'
$ws.Range('E20').Value = '  +8.74%  '
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').Value = '67.79'
$ws.Range('E23').Value = '  +1.80%  '
$ws.Range('D24').Value = '236.71'
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('E25').Value = '  +12.51%  '
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('D28').Value = '24.66'
$ws.Range('E28').Value = '  +3.98%  '
$ws.Range('D29').Value = '168.18'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').Value = '2.08'
$ws.Range('E30').Value = '  -8.87%  '
$ws.Range('D31').Value = '33.76'
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('D32').Value = '9.16'
$ws.Range('E32').Value = '  +1.43%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('E35').Value = '  +4.47%  '
$ws.Range('E36').Value = '  +3.74%  '
$ws.Range('D37').Value = '16.95'
$ws.Range('E37').Value = '  +6.57%  '
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('E39').Value = '  +3.58%  '
$ws.Range('E40').Value = '  +4.99%  '
$ws.Range('D41').Value = '2.81'
$ws.Range('E41').Value = '  +1.20%  '
$ws.Range('D42').Value = '0.109'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('E43').Value = '  -2.92%  '
$ws.Range('D44').Value = '2.002.26'
$ws.Range('E44').Value = '  +2.84%  '
$ws.Range('D45').Value = '0.0286'
$ws.Range('E45').Value = '  +2.59%  '
$ws.Range('E46').Value = '  +5.59%  '
$ws.Range('D47').Value = '17.73'
$ws.Range('E47').Value = '  +2.20%  '
$ws.Range('E48').Value = '  +2.91%  '
$ws.Range('D49').Value = '55.13'
$ws.Range('E49').Value = '  +5.96%  '
$ws.Range('D50').Value = '2.527.96'
$ws.Range('E50').Value = '  +1.74%  '
$ws.Range('E51').Value = '  +2.09%  '
